$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.598.63"
$ws.Range("E2").Value = "  -0.18%  "

$ws.Range("D3").Value = "3.506.18"
$ws.Range("E3").Value = "  -0.37%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "195.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.72%  "

$ws.Range("E7").Value = "  -1.25%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.200"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.88%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.644"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.77%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.09"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.76%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000298"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.39%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.59%  "

$ws.Range("D14").Value = "4.060.23"
$ws.Range("E14").Value = "  -0.56%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "598.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.30%  "

$ws.Range("D16").Value = "69.706.50"
$ws.Range("E16").Value = "  -0.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.06%  "

$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.529.16"
$ws.Range("E19").Value = "  +0.41%  "

$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.123"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.981"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.81%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.96%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "101.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.85%  "

$ws.Range("E25").Value = "  -2.04%  "

$ws.Range("E26").Value = "  +0.32%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.52%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.57%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.50%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.73%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.114"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.56%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.07%  "

$ws.Range("D36").Value = "3.729.20"
$ws.Range("E36").Value = "  +1.60%  "

$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0808"
$ws.Range("E37").Value = "  +2.51%  "

$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.95%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.388"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "492.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.60%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.132"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.72%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0448"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.52%  "

$ws.Range("E45").Value = "  -4.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.139"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.84%  "

$ws.Range("E47").Value = "  -2.90%  "

$ws.Range("E48").Value = "  +0.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.17%  "

$ws.Range("E50").Value = "  +0.55%  "

$ws.Range("B51").Value = "OceanProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.89%  "
